{"js": "// Replace the multiplication-problem strings in the practice table.\n// Each cell holds a single run of text like \"840\u00d78=\" \u2014 we search the\n// whole document body for each old value and replace it with its new\n// value. All 25 old values are unique in the document, so a plain\n// text search safely targets the correct cell each time.\nconst replacements = [\n  [\"840\u00d78=\", \"662\u00d72=\"],\n  [\"630\u00d79=\", \"251\u00d76=\"],\n  [\"605\u00d77=\", \"409\u00d76=\"],\n  [\"953\u00d77=\", \"477\u00d72=\"],\n  [\"386\u00d79=\", \"156\u00d77=\"],\n  [\"939\u00d74=\", \"787\u00d79=\"],\n  [\"185\u00d79=\", \"301\u00d74=\"],\n  [\"542\u00d75=\", \"140\u00d74=\"],\n  [\"792\u00d77=\", \"468\u00d75=\"],\n  [\"392\u00d77=\", \"286\u00d74=\"],\n  [\"675\u00d79=\", \"151\u00d77=\"],\n  [\"748\u00d75=\", \"186\u00d74=\"],\n  [\"756\u00d73=\", \"790\u00d73=\"],\n  [\"115\u00d72=\", \"596\u00d73=\"],\n  [\"329\u00d74=\", \"737\u00d73=\"],\n  [\"471\u00d78=\", \"385\u00d79=\"],\n  [\"949\u00d75=\", \"589\u00d73=\"],\n  [\"147\u00d73=\", \"471\u00d77=\"],\n  [\"941\u00d79=\", \"679\u00d77=\"],\n  [\"427\u00d78=\", \"478\u00d77=\"],\n  [\"695\u00d74=\", \"957\u00d78=\"],\n  [\"880\u00d76=\", \"169\u00d74=\"],\n  [\"758\u00d77=\", \"758\u00d73=\"],\n  [\"151\u00d78=\", \"699\u00d79=\"],\n  [\"642\u00d77=\", \"598\u00d73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each \"AAA\u00d7B=\" multiplication prompt in the practice table to its\n# new value. Every old value is unique within the document, so a plain\n# Find/Replace (no wildcards) targets exactly one cell each time.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"840\u00d78=\", \"662\u00d72=\"),\n  @(\"630\u00d79=\", \"251\u00d76=\"),\n  @(\"605\u00d77=\", \"409\u00d76=\"),\n  @(\"953\u00d77=\", \"477\u00d72=\"),\n  @(\"386\u00d79=\", \"156\u00d77=\"),\n  @(\"939\u00d74=\", \"787\u00d79=\"),\n  @(\"185\u00d79=\", \"301\u00d74=\"),\n  @(\"542\u00d75=\", \"140\u00d74=\"),\n  @(\"792\u00d77=\", \"468\u00d75=\"),\n  @(\"392\u00d77=\", \"286\u00d74=\"),\n  @(\"675\u00d79=\", \"151\u00d77=\"),\n  @(\"748\u00d75=\", \"186\u00d74=\"),\n  @(\"756\u00d73=\", \"790\u00d73=\"),\n  @(\"115\u00d72=\", \"596\u00d73=\"),\n  @(\"329\u00d74=\", \"737\u00d73=\"),\n  @(\"471\u00d78=\", \"385\u00d79=\"),\n  @(\"949\u00d75=\", \"589\u00d73=\"),\n  @(\"147\u00d73=\", \"471\u00d77=\"),\n  @(\"941\u00d79=\", \"679\u00d77=\"),\n  @(\"427\u00d78=\", \"478\u00d77=\"),\n  @(\"695\u00d74=\", \"957\u00d78=\"),\n  @(\"880\u00d76=\", \"169\u00d74=\"),\n  @(\"758\u00d77=\", \"758\u00d73=\"),\n  @(\"151\u00d78=\", \"699\u00d79=\"),\n  @(\"642\u00d77=\", \"598\u00d73=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
